$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) values stay as plain text (they are not real numbers,
# e.g. "38.183.81" or "0.999"), matching the source data which stores prices as text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.183.81'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.114.19'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.23'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.29'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0782'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.416.32'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.71'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.74'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.789'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.29'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.109.84'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.113.49'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.18'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.93'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0829'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.10'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.48'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.142'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.05'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.56'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.69'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.62'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0630'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.63'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.50'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.51'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0997'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.51'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0217'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.468.05'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.16'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.02'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.07'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.34'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.311.72'

# Update coin name / link / volume text cells
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("E3").Value = '  +3.01%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("E6").Value = '  +1.30%  '
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +2.89%  '
$ws.Range("E10").Value = '  +3.39%  '
$ws.Range("E11").Value = '  +2.16%  '
$ws.Range("E12").Value = '  +2.57%  '
$ws.Range("E13").Value = '  +3.12%  '
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("E16").Value = '  +3.13%  '
$ws.Range("E17").Value = '  +2.75%  '
$ws.Range("E18").Value = '  +2.64%  '
$ws.Range("E19").Value = '  -2.87%  '
$ws.Range("E20").Value = '  +2.86%  '
$ws.Range("E21").Value = '  +2.69%  '
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("E26").Value = '  +1.37%  '
$ws.Range("E27").Value = '  +13.00%  '
$ws.Range("E28").Value = '  +3.30%  '
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("E32").Value = '  +5.52%  '
$ws.Range("E33").Value = '  +4.43%  '
$ws.Range("E34").Value = '  +2.41%  '
$ws.Range("E35").Value = '  +1.29%  '
$ws.Range("E36").Value = '  +6.79%  '
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("E39").Value = '  -3.51%  '
$ws.Range("E40").Value = '  +7.40%  '
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("E43").Value = '  +3.36%  '
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("E46").Value = '  -6.48%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E47").Value = '  +6.21%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E48").Value = '  +5.16%  '
$ws.Range("E49").Value = '  +3.61%  '
$ws.Range("E50").Value = '  +2.22%  '
$ws.Range("E51").Value = '  +3.03%  '

Write-Host "Applied cryptos update"
